# Data-cleaning pass over the purchases table:
#  - normalise category labels (accents/typos/case) in column E
#  - fix a handful of mis-typed product names in column B
#  - correct several purchase dates in column F
#  - reselect the full table range and tidy the Categorie column width
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Produit (B): fix mis-typed / inconsistent product names ---
$ws.Range("B7").Value = "pates"
$ws.Range("B8").Value = "pates"
$ws.Range("B13").Value = "oeuf"
$ws.Range("B18").Value = "banane"
$ws.Range("B20").Value = "pomme"
$ws.Range("B22").Value = "coca cola"
$ws.Range("B24").Value = "eau  minerale"
$ws.Range("B26").Value = "cafe"
$ws.Range("B27").Value = "cafe"
$ws.Range("B28").Value = "the"
$ws.Range("B29").Value = "the"
$ws.Range("B30").Value = "choco lait"
$ws.Range("B35").Value = "pates"
$ws.Range("B44").Value = "eau  minerale"
$ws.Range("B47").Value = "corn flakes"
$ws.Range("B48").Value = "corn flakes"

# --- Categorie (E): normalise category labels ---
$ws.Range("E2").Value = "Boulangerie"
$ws.Range("E3").Value = "Laitage"
$ws.Range("E4").Value = "Laitage"
$ws.Range("E5").Value = "Fruits & Légumes"
$ws.Range("E6").Value = "Fruits & Légumes"
$ws.Range("E7").Value = "Epicerie"
$ws.Range("E8").Value = "Epicerie"
$ws.Range("E9").Value = "Epicerie"
$ws.Range("E10").Value = "Epicerie"
$ws.Range("E11").Value = "Laitage"
$ws.Range("E12").Value = "Laitage"
$ws.Range("E13").Value = "Œufs & Ovoproduits"
$ws.Range("E14").Value = "Œufs & Ovoproduits"
$ws.Range("E15").Value = "Boucherie"
$ws.Range("E16").Value = "Poissonnerie"
$ws.Range("E17").Value = "Fruits & Légumes"
$ws.Range("E18").Value = "Fruits & Légumes"
$ws.Range("E19").Value = "Fruits & Légumes"
$ws.Range("E20").Value = "Fruits & Légumes"
$ws.Range("E21").Value = "Fruits & Légumes"
$ws.Range("E22").Value = "Boissons"
$ws.Range("E23").Value = "Boissons"
$ws.Range("E24").Value = "Boissons"
$ws.Range("E25").Value = "Boissons"
$ws.Range("E26").Value = "Epicerie"
$ws.Range("E27").Value = "Epicerie"
$ws.Range("E28").Value = "Epicerie"
$ws.Range("E29").Value = "Epicerie"
$ws.Range("E30").Value = "Epicerie"
$ws.Range("E31").Value = "Epicerie"
$ws.Range("E32").Value = "Laitage"
$ws.Range("E33").Value = "Laitage"
$ws.Range("E34").Value = "Boulangerie"
$ws.Range("E35").Value = "Epicerie"
$ws.Range("E36").Value = "Epicerie"
$ws.Range("E37").Value = "Fruits & Légumes"
$ws.Range("E38").Value = "Cremerie"
$ws.Range("E39").Value = "Cremerie"
$ws.Range("E40").Value = "Charcuterie"
$ws.Range("E41").Value = "Charcuterie"
$ws.Range("E42").Value = "Laitage"
$ws.Range("E43").Value = "Laitage"
$ws.Range("E44").Value = "Boissons"
$ws.Range("E45").Value = "Epicerie"
$ws.Range("E46").Value = "Epicerie"
$ws.Range("E47").Value = "Epicerie"
$ws.Range("E48").Value = "Epicerie"
$ws.Range("E49").Value = "Laitage"

# --- Date (F): correct purchase dates ---
$ws.Range("F3").Value = 45666
$ws.Range("F6").Value = 45697
$ws.Range("F8").Value = 45697
$ws.Range("F12").Value = 45725
$ws.Range("F14").Value = 45756
$ws.Range("F18").Value = 45786
$ws.Range("F23").Value = 45817
$ws.Range("F25").Value = 45817
$ws.Range("F27").Value = 45847
$ws.Range("F41").Value = 45939
$ws.Range("F43").Value = 45666
$ws.Range("F46").Value = 45970
$ws.Range("F49").Value = 46000

# --- Reselect the full table range ---
$ws.Range("A1:F49").Select()

# --- Tidy the Categorie column width for the relabeled content ---
$ws.Columns.Item(5).ColumnWidth = 16.8
